$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.747.56"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "2.650.89"
$ws.Range("E3").Value = "  +3.62%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "513.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.566"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.21%  "
$ws.Range("D9").Value = "2.678.43"
$ws.Range("E9").Value = "  +4.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.105"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.335"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("D14").Value = "3.123.74"
$ws.Range("E14").Value = "  +3.84%  "
$ws.Range("D15").Value = "58.785.06"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.33%  "
$ws.Range("E17").Value = "  +1.96%  "
$ws.Range("D18").Value = "2.681.83"
$ws.Range("E18").Value = "  +4.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "342.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.419"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.83%  "
$ws.Range("D26").Value = "2.770.14"
$ws.Range("E26").Value = "  +3.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.991"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.160"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.63%  "
$ws.Range("D29").Value = "0.0₃0809"
$ws.Range("E29").Value = "  +4.00%  "
$ws.Range("E30").Value = "  +4.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.64%  "
$ws.Range("E34").Value = "  +2.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.49"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.01"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.69%  "
$ws.Range("E38").Value = "  +3.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.853"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.50"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.67"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.90%  "
$ws.Range("E42").Value = "  +1.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "280.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.618"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.993"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0978"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0532"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0230"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.70"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.95%  "
